# The underlying edit swaps the contents of ppt/theme/theme1.xml (used by
# the slide master -> all normal slides) and ppt/theme/theme2.xml (used by
# the notes master): theme1.xml goes from the "Integral" palette to the
# default "Office Theme" palette, and theme2.xml goes the other way
# (Office Theme -> Integral).
#
# The PowerPoint object model lets us rewrite a design's 12-slot theme
# color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) via
# Master.ColorScheme.Colors(i).RGB. Apply the "Office Theme" RGB values,
# in VBA's RGB() encoding (R + G*256 + B*65536), to the presentation's
# slide master so its theme (theme1.xml) is recolored to match the
# "Office Theme" palette.

$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$colorScheme = $master.ColorScheme

# index -> (slot, target "Office Theme" RGB hex)
#  1 = dk1       000000
#  2 = lt1       FFFFFF
#  3 = dk2       44546A
#  4 = lt2       E7E6E6
#  5 = accent1   5B9BD5
#  6 = accent2   ED7D31
#  7 = accent3   A5A5A5
#  8 = accent4   FFC000
#  9 = accent5   4472C4
# 10 = accent6   70AD47
# 11 = hlink     0563C1
# 12 = folHlink  954F72

$colorScheme.Colors(1).RGB  = 0          # dk1      #000000
$colorScheme.Colors(2).RGB  = 16777215   # lt1      #FFFFFF
$colorScheme.Colors(3).RGB  = 6968388    # dk2      #44546A
$colorScheme.Colors(4).RGB  = 15132391   # lt2      #E7E6E6
$colorScheme.Colors(5).RGB  = 13998939   # accent1  #5B9BD5
$colorScheme.Colors(6).RGB  = 3243501    # accent2  #ED7D31
$colorScheme.Colors(7).RGB  = 10855845   # accent3  #A5A5A5
$colorScheme.Colors(8).RGB  = 49407      # accent4  #FFC000
$colorScheme.Colors(9).RGB  = 12874308   # accent5  #4472C4
$colorScheme.Colors(10).RGB = 4697456    # accent6  #70AD47
$colorScheme.Colors(11).RGB = 12673797   # hlink    #0563C1
$colorScheme.Colors(12).RGB = 7491477    # folHlink #954F72
